# Slide 25 ("Follow Sets from CPRL: Example 1"), body placeholder
# (Shapes.Item(2), "Rectangle 3"), third paragraph:
#
#   "subprogramDecls = subprogramDecl { subprogramDecl } ."
#
# needs to become (matching the rest of the deck's updated grammar style):
#
#   "subprogramDecls = { subprogramDecl } ."
#
# split across three runs: "subprogramDecls ", "= { ", "subprogramDecl } ."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(25)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(3)

$oldMiddle = "= subprogramDecl { "
$newMiddle = "= { "

$startIdx = $para.Text.IndexOf($oldMiddle)
if ($startIdx -lt 0) {
    throw "Expected substring not found in paragraph: $($para.Text)"
}

# Characters() is 1-indexed, IndexOf() is 0-indexed.
$midRange = $para.Characters($startIdx + 1, $oldMiddle.Length)
$midRange.Text = $newMiddle

Write-Output "Updated paragraph text: $($para.Text)"
